$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 originally held: A=player-gamelog URL, B=date, C=week, D=fantasy
# points allowed, E=team, F=(blank), G=opponent, H=result, I=(blank),
# J=rushes/receptions stat, K=week (dup), L=yards-per-stat, M=fumbles -
# all stored as plain text, even the numeric-looking ones.
#
# The edit drops the URL column, adds three new leading identity columns
# (last name, first name, position), shifts the remaining original values
# two columns to the right (still as literal text), and appends a new
# numeric fantasy-points-per-stat value at the end.

function Set-TextValue($addr, $text) {
    # Force the literal string into the cell (instead of letting Excel
    # auto-parse look-alike numbers/dates), then drop back to the default
    # "Normal" style so no stray number-format styling is left behind.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "A1" "McNichols"
Set-TextValue "B1" "Jeremy"
Set-TextValue "C1" "RB"
Set-TextValue "D1" "2018-09-30"
Set-TextValue "E1" "4"
Set-TextValue "F1" "22.278"
Set-TextValue "G1" "IND"
Set-TextValue "H1" ""
Set-TextValue "I1" "HOU"
Set-TextValue "J1" "L 34-37"
Set-TextValue "K1" ""
Set-TextValue "L1" "2"
Set-TextValue "M1" "4"
Set-TextValue "N1" "2.00"
Set-TextValue "O1" "0"

$ws.Range("P1").Value = 0.4
